# Update tab names in template4.xlsx
#
# 1) Rename sheet "ETPT_ATT_JUR" -> "ETPT_ATTJ" and unhide it.
# 2) Rename sheet "ETPT_ATT_JUR_DDG" -> "ETPT_ATTJ_DDG".
#    (Excel automatically rewrites every formula referencing these sheets.)
# 3) Update the selections left on a few sheets to cell I5.
# 4) Make "ETPT A-JUST" the active sheet/tab (instead of "ETPT_ATT_JUR_DDG").

$wb = $excel.ActiveWorkbook

# --- Rename the two "ATT_JUR" sheets -------------------------------------
$wsAttJur = $wb.Worksheets.Item("ETPT_ATT_JUR")
$wsAttJur.Name = "ETPT_ATTJ"
$wsAttJur.Visible = -1

$wsAttJurDdg = $wb.Worksheets.Item("ETPT_ATT_JUR_DDG")
$wsAttJurDdg.Name = "ETPT_ATTJ_DDG"

# --- Update leftover cell selections on a few sheets ----------------------
$wsAttJur.Range("I5").Select()

$wsTjDdg = $wb.Worksheets.Item("ETPT_TJ_DDG")
$wsTjDdg.Range("I5").Select()

$wsAttJurDdg.Range("I5").Select()

# --- Make "ETPT A-JUST" the active sheet ----------------------------------
$wsActive = $wb.Worksheets.Item("ETPT A-JUST")
$wsActive.Activate()
